$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5517241379310345
$ws.Range("C2").Value = 0.8421052631578947
$ws.Range("D2").Value = 0.6666666666666666
$ws.Range("B3").Value = 0.8285714285714286
$ws.Range("C3").Value = 0.5272727272727272
$ws.Range("D3").Value = 0.6444444444444444
$ws.Range("B4").Value = 0.6559139784946236
$ws.Range("C4").Value = 0.6559139784946236
$ws.Range("D4").Value = 0.6559139784946236
$ws.Range("E4").Value = 0.6559139784946236
$ws.Range("B5").Value = 0.6901477832512315
$ws.Range("C5").Value = 0.684688995215311
$ws.Range("D5").Value = 0.6555555555555554
$ws.Range("B6").Value = 0.715451030245246
$ws.Range("C6").Value = 0.6559139784946236
$ws.Range("D6").Value = 0.6535244922341696
$ws.Range("B7").Value = 0.543859649122807
$ws.Range("C7").Value = 0.8157894736842105
$ws.Range("D7").Value = 0.6526315789473685
$ws.Range("B8").Value = 0.8055555555555556
$ws.Range("C8").Value = 0.5272727272727272
$ws.Range("D8").Value = 0.6373626373626373
$ws.Range("B9").Value = 0.6451612903225806
$ws.Range("C9").Value = 0.6451612903225806
$ws.Range("D9").Value = 0.6451612903225806
$ws.Range("E9").Value = 0.6451612903225806
$ws.Range("B10").Value = 0.6747076023391814
$ws.Range("C10").Value = 0.6715311004784689
$ws.Range("D10").Value = 0.6449971081550029
$ws.Range("B11").Value = 0.698626045400239
$ws.Range("C11").Value = 0.6451612903225806
$ws.Range("D11").Value = 0.643601559730592
$ws.Range("B12").Value = 0.5178571428571429
$ws.Range("C12").Value = 0.7631578947368421
$ws.Range("D12").Value = 0.6170212765957447
$ws.Range("B13").Value = 0.7567567567567568
$ws.Range("C13").Value = 0.509090909090909
$ws.Range("D13").Value = 0.608695652173913
$ws.Range("B14").Value = 0.6129032258064516
$ws.Range("C14").Value = 0.6129032258064516
$ws.Range("D14").Value = 0.6129032258064516
$ws.Range("E14").Value = 0.6129032258064516
$ws.Range("B15").Value = 0.6373069498069499
$ws.Range("C15").Value = 0.6361244019138756
$ws.Range("D15").Value = 0.6128584643848288
$ws.Range("B16").Value = 0.659141860754764
$ws.Range("C16").Value = 0.6129032258064516
$ws.Range("D16").Value = 0.6120975202172421
$ws.Range("B17").Value = 0.5428571428571428
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0.7037037037037037
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0.4181818181818182
$ws.Range("D18").Value = 0.5897435897435896
$ws.Range("B19").Value = 0.6559139784946236
$ws.Range("C19").Value = 0.6559139784946236
$ws.Range("D19").Value = 0.6559139784946236
$ws.Range("E19").Value = 0.6559139784946236
$ws.Range("B20").Value = 0.7714285714285714
$ws.Range("C20").Value = 0.7090909090909091
$ws.Range("D20").Value = 0.6467236467236467
$ws.Range("B21").Value = 0.8132104454685101
$ws.Range("C21").Value = 0.6559139784946236
$ws.Range("D21").Value = 0.6363079373832061
$ws.Range("B22").Value = 0.625
$ws.Range("C22").Value = 0.7894736842105263
$ws.Range("D22").Value = 0.6976744186046512
$ws.Range("B23").Value = 0.8222222222222222
$ws.Range("C23").Value = 0.6727272727272727
$ws.Range("D23").Value = 0.7400000000000001
$ws.Range("B24").Value = 0.7204301075268817
$ws.Range("C24").Value = 0.7204301075268817
$ws.Range("D24").Value = 0.7204301075268817
$ws.Range("E24").Value = 0.7204301075268817
$ws.Range("B25").Value = 0.7236111111111111
$ws.Range("C25").Value = 0.7311004784688995
$ws.Range("D25").Value = 0.7188372093023256
$ws.Range("B26").Value = 0.7416367980884111
$ws.Range("C26").Value = 0.7204301075268817
$ws.Range("D26").Value = 0.7227056764191048
